$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("G9").Value = 2.1
$ws.Range("I9").Value = 3.2
$ws.Range("K9").Value = 2.38
$ws.Range("X9").Value = 12
$ws.Range("AC9").Value = 15
$ws.Range("AD9").Value = 7
$ws.Range("AE9").Value = 11
$ws.Range("AI9").Value = 12
$ws.Range("AK9").Value = 23
$ws.Range("AL9").Value = 26
$ws.Range("AX9").Value = 17

# Row 16
$ws.Range("G16").Value = 2.38
$ws.Range("Y16").Value = 9.5
$ws.Range("AK16").Value = 26
$ws.Range("AN16").Value = 4.33
$ws.Range("AW16").Value = 5

# Row 24
$ws.Range("H24").Value = 3.75
$ws.Range("I24").Value = 4.55
$ws.Range("J24").Value = 2.18
$ws.Range("L24").Value = 4.85
$ws.Range("Q24").Value = 1.72
$ws.Range("R24").Value = 2.05
$ws.Range("V24").Value = 1.98
$ws.Range("W24").Value = 7.6
$ws.Range("X24").Value = 8.25
$ws.Range("Y24").Value = 8
$ws.Range("AB24").Value = 23
$ws.Range("AD24").Value = 7.5
$ws.Range("AE24").Value = 15
$ws.Range("AH24").Value = 28
$ws.Range("AJ24").Value = 80
$ws.Range("AK24").Value = 45
$ws.Range("AM24").Value = 450
$ws.Range("AW24").Value = 6.4
$ws.Range("AX24").Value = 26

# Row 26
$ws.Range("G26").Value = 2.95
$ws.Range("H26").Value = 2.72
$ws.Range("I26").Value = 2.6
$ws.Range("J26").Value = 3.65
$ws.Range("K26").Value = 1.83
$ws.Range("L26").Value = 3.35
$ws.Range("N26").Value = 5.7
$ws.Range("U26").Value = 2.02
$ws.Range("V26").Value = 1.62
$ws.Range("W26").Value = 6.9
$ws.Range("X26").Value = 14
$ws.Range("Y26").Value = 11.25
$ws.Range("Z26").Value = 40
$ws.Range("AA26").Value = 32
$ws.Range("AB26").Value = 50
$ws.Range("AC26").Value = 5.8
$ws.Range("AD26").Value = 5.5
$ws.Range("AE26").Value = 17.5
$ws.Range("AG26").Value = 6.1
$ws.Range("AH26").Value = 11.5
$ws.Range("AI26").Value = 10.5
$ws.Range("AJ26").Value = 30
$ws.Range("AK26").Value = 29
$ws.Range("AN26").Value = 4.6
$ws.Range("AO26").Value = 17
$ws.Range("AP26").Value = 28
$ws.Range("AQ26").Value = 90
$ws.Range("AR26").Value = 150
$ws.Range("AS26").Value = 450
$ws.Range("AU26").Value = 7.5
$ws.Range("AV26").Value = 90
$ws.Range("AW26").Value = 4.25
$ws.Range("AX26").Value = 15
$ws.Range("AY26").Value = 27
$ws.Range("AZ26").Value = 75
